$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Name = "ProductName"
